# Archivio fine settimana lavorativa
# Update simulation results with refreshed solver output values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = [double]"389.4877318255876"
$ws.Range("C5").Value = [double]"79.99999978583787"
$ws.Range("D5").Value = [double]"38.74576745149522"
$ws.Range("E5").Value = [double]"131.7628756641556"
$ws.Range("B6").Value = [double]"389.4877318255876"
$ws.Range("C6").Value = [double]"119.9999996787568"
$ws.Range("D6").Value = [double]"58.11865117724281"
$ws.Range("E6").Value = [double]"197.6443134962334"
$ws.Range("B7").Value = [double]"390.9920215025131"
$ws.Range("C7").Value = [double]"39.99999980483442"
$ws.Range("D7").Value = [double]"19.37288266202678"
$ws.Range("E7").Value = [double]"65.62796698143455"
$ws.Range("B8").Value = [double]"390.9920215025131"
$ws.Range("C8").Value = [double]"59.99999970725163"
$ws.Range("D8").Value = [double]"29.05932399304016"
$ws.Range("E8").Value = [double]"98.44195047215182"
$ws.Range("B9").Value = [double]"390.3614997407959"
$ws.Range("C9").Value = [double]"39.99999920291975"
$ws.Range("D9").Value = [double]"-1.192816870968727E-06"
$ws.Range("E9").Value = [double]"59.16057378286122"
$ws.Range("B10").Value = [double]"383.6890695056338"
$ws.Range("C10").Value = [double]"199.9999914477247"
$ws.Range("D10").Value = [double]"96.86440775791222"
$ws.Range("E10").Value = [double]"334.3854757602283"
$ws.Range("B11").Value = [double]"390.6949334566526"
$ws.Range("C11").Value = [double]"149.9999991032052"
$ws.Range("D11").Value = [double]"72.64830986940409"
$ws.Range("E11").Value = [double]"246.2920163145269"
$ws.Range("B15").Value = [double]"1492.664422274971"
$ws.Range("C15").Value = [double]"199.9999998133614"
$ws.Range("D15").Value = [double]"-2.324408141021195E-07"
$ws.Range("E15").Value = [double]"133.9885890149003"
$ws.Range("B16").Value = [double]"1492.664422274971"
$ws.Range("C16").Value = [double]"199.9999998133614"
$ws.Range("D16").Value = [double]"-2.324408061440409E-07"
$ws.Range("E16").Value = [double]"133.9885890149003"
$ws.Range("B17").Value = [double]"390.9920572034896"
$ws.Range("C17").Value = [double]"49.99999905139003"
$ws.Range("D17").Value = [double]"-1.454926377505217E-06"
$ws.Range("E17").Value = [double]"127.8798331839534"
$ws.Range("B18").Value = [double]"586.0423178031964"
$ws.Range("C18").Value = [double]"29.99999939523494"
$ws.Range("D18").Value = [double]"-8.783536705436746E-07"
$ws.Range("E18").Value = [double]"51.19084148677153"
$ws.Range("B22").Value = [double]"1492.664422274971"
$ws.Range("C22").Value = [double]"9.999999990668071"
$ws.Range("D22").Value = [double]"-1.162204035409786E-08"
$ws.Range("E22").Value = [double]"6.699429450745013"
$ws.Range("B23").Value = [double]"586.4879829573833"
$ws.Range("C23").Value = [double]"-9.999999810278005"
$ws.Range("D23").Value = [double]"2.909852749013453E-07"
$ws.Range("E23").Value = [double]"17.05064741453582"
$ws.Range("B24").Value = [double]"586.0423178039709"
$ws.Range("C24").Value = [double]"49.99999899205824"
$ws.Range("D24").Value = [double]"-1.463922787138472E-06"
$ws.Range("E24").Value = [double]"85.31806914450645"
$ws.Range("B25").Value = [double]"602.9445885376979"
$ws.Range("C25").Value = [double]"99.99999993000156"
$ws.Range("D25").Value = [double]"-4.56406350792804E-08"
$ws.Range("E25").Value = [double]"165.8527198536265"
$ws.Range("B29").Value = [double]"746.3322131690232"
$ws.Range("C29").Value = [double]"-14.99999998600211"
$ws.Range("D29").Value = [double]"1.743306001600331E-08"
$ws.Range("E29").Value = [double]"20.0982882975266"
$ws.Range("B30").Value = [double]"746.3322131690232"
$ws.Range("C30").Value = [double]"-14.99999998600211"
$ws.Range("D30").Value = [double]"1.743306001600331E-08"
$ws.Range("E30").Value = [double]"20.0982882975266"
$ws.Range("B31").Value = [double]"586.4879829573833"
$ws.Range("C31").Value = [double]"-29.99999943083401"
$ws.Range("D31").Value = [double]"8.729558271625138E-07"
$ws.Range("E31").Value = [double]"51.15194224360746"
$ws.Range("B32").Value = [double]"586.0423178031964"
$ws.Range("C32").Value = [double]"-49.99999899205825"
$ws.Range("D32").Value = [double]"1.463922786683725E-06"
$ws.Range("E32").Value = [double]"85.31806914461924"
$ws.Range("B33").Value = [double]"602.9445885376979"
$ws.Range("C33").Value = [double]"-39.99999997200063"
$ws.Range("D33").Value = [double]"1.825625402318565E-08"
$ws.Range("E33").Value = [double]"66.3410879414506"
$ws.Range("B37").Value = [double]"2029.411106824899"
$ws.Range("C37").Value = [double]"-139.9999999227312"
$ws.Range("D37").Value = [double]"-90.43106628120154"
$ws.Range("E37").Value = [double]"47.41525482618598"
$ws.Range("B38").Value = [double]"2029.411106824899"
$ws.Range("C38").Value = [double]"-139.999999916552"
$ws.Range("D38").Value = [double]"-79.34165617088203"
$ws.Range("E38").Value = [double]"45.78024603907611"
$ws.Range("B42").Value = [double]"586.0423178039709"
$ws.Range("C42").Value = [double]"-9.999999798411647"
$ws.Range("D42").Value = [double]"2.9278455748738E-07"
$ws.Range("E42").Value = [double]"17.06361382890129"
$ws.Range("B49").Value = [double]"19981.41654245573"
$ws.Range("C49").Value = [double]"1990.219221749981"
$ws.Range("D49").Value = [double]"291.0978451837209"
$ws.Range("E49").Value = [double]"290.6670371781187"
$ws.Range("F49").Value = [double]"87.95996414089622"
$ws.Range("G49").Value = [double]"55.47160769024637"
$ws.Range("H49").Value = [double]"8.786685012030956"
$ws.Range("I49").Value = [double]"85.84249010011035"
$ws.Range("B50").Value = [double]"19981.41654245573"
$ws.Range("C50").Value = [double]"1990.219221749981"
$ws.Range("D50").Value = [double]"291.0978451837209"
$ws.Range("E50").Value = [double]"290.6670371781187"
$ws.Range("F50").Value = [double]"87.95996414089622"
$ws.Range("G50").Value = [double]"55.47160769024637"
$ws.Range("H50").Value = [double]"8.786685012030956"
$ws.Range("I50").Value = [double]"85.84249010011035"
$ws.Range("B51").Value = [double]"1990.219221749981"
$ws.Range("C51").Value = [double]"389.4877318255876"
$ws.Range("D51").Value = [double]"201.3340750800523"
$ws.Range("E51").Value = [double]"199.9999974099957"
$ws.Range("F51").Value = [double]"110.9372349975222"
$ws.Range("G51").Value = [double]"96.8644196053078"
$ws.Range("H51").Value = [double]"66.6852909752578"
$ws.Range("I51").Value = [double]"329.4071870503547"
$ws.Range("B52").Value = [double]"19991.44039087657"
$ws.Range("C52").Value = [double]"390.9920215025131"
$ws.Range("D52").Value = [double]"216.7774766492043"
$ws.Range("E52").Value = [double]"215.2759859109946"
$ws.Range("F52").Value = [double]"114.0867592606075"
$ws.Range("G52").Value = [double]"98.96249369948049"
$ws.Range("H52").Value = [double]"7.074625567814823"
$ws.Range("I52").Value = [double]"349.8691375430906"
$ws.Range("B53").Value = [double]"19996.52758431912"
$ws.Range("C53").Value = [double]"392.630725570652"
$ws.Range("D53").Value = [double]"311.2948968348108"
$ws.Range("E53").Value = [double]"309.4282645542194"
$ws.Range("F53").Value = [double]"143.7099561083788"
$ws.Range("G53").Value = [double]"122.7976353683374"
$ws.Range("H53").Value = [double]"9.899474712498384"
$ws.Range("I53").Value = [double]"489.5308855231982"
$ws.Range("B54").Value = [double]"19996.52758431912"
$ws.Range("C54").Value = [double]"2009.81528168955"
$ws.Range("D54").Value = [double]"-217.7305242072229"
$ws.Range("E54").Value = [double]"-218.0358415622719"
$ws.Range("F54").Value = [double]"-135.0188439554549"
$ws.Range("G54").Value = [double]"-167.2854910187111"
$ws.Range("H54").Value = [double]"7.397041424311022"
$ws.Range("I54").Value = [double]"78.94523177262131"
$ws.Range("B58").Value = [double]"1990.219221749981"
$ws.Range("C58").Value = [double]"1492.664422274971"
$ws.Range("D58").Value = [double]"199.994363128181"
$ws.Range("E58").Value = [double]"200.003438768862"
$ws.Range("F58").Value = [double]"0.006348268891982003"
$ws.Range("G58").Value = [double]"0.008778492124223818"
$ws.Range("H58").Value = [double]"58.01514330357119"
$ws.Range("I58").Value = [double]"133.9902314568286"
$ws.Range("B59").Value = [double]"1990.219221749981"
$ws.Range("C59").Value = [double]"1492.664422274971"
$ws.Range("D59").Value = [double]"179.9941523191701"
$ws.Range("E59").Value = [double]"180.0022935533734"
$ws.Range("F59").Value = [double]"0.002939366966476428"
$ws.Range("G59").Value = [double]"0.004415428849018781"
$ws.Range("H59").Value = [double]"52.21562209227531"
$ws.Range("B60").Value = [double]"390.9920215025131"
$ws.Range("C60").Value = [double]"586.4879829573833"
$ws.Range("D60").Value = [double]"75.2182949705508"
$ws.Range("E60").Value = [double]"75.21943558418545"
$ws.Range("F60").Value = [double]"50.50333516189703"
$ws.Range("G60").Value = [double]"50.47960675645903"
$ws.Range("H60").Value = [double]"133.8033844980201"
$ws.Range("I60").Value = [double]"154.404914706582"
$ws.Range("B61").Value = [double]"390.6949334566526"
$ws.Range("C61").Value = [double]"586.0423178039709"
$ws.Range("D61").Value = [double]"-25.15875725152281"
$ws.Range("E61").Value = [double]"-25.15971686606806"
$ws.Range("F61").Value = [double]"-50.44664287785054"
$ws.Range("G61").Value = [double]"-50.44165817080625"
$ws.Range("H61").Value = [double]"83.29330985690835"
$ws.Range("I61").Value = [double]"96.17403605037291"
$ws.Range("B62").Value = [double]"390.6949334566526"
$ws.Range("C62").Value = [double]"586.0423178031964"
$ws.Range("D62").Value = [double]"-20.0014096484914"
$ws.Range("E62").Value = [double]"-20.00245079530523"
$ws.Range("F62").Value = [double]"-0.002231392821311033"
$ws.Range("G62").Value = [double]"0.0007522937014116734"
$ws.Range("H62").Value = [double]"29.55968424833608"
$ws.Range("I62").Value = [double]"34.12833818153045"
$ws.Range("B63").Value = [double]"2009.81528168955"
$ws.Range("C63").Value = [double]"602.9445885376979"
$ws.Range("D63").Value = [double]"59.99987326924993"
$ws.Range("E63").Value = [double]"60.00021625807812"
$ws.Range("F63").Value = [double]"-0.004836555941451451"
$ws.Range("G63").Value = [double]"-0.001124051604098369"
$ws.Range("H63").Value = [double]"17.23624520930266"
$ws.Range("I63").Value = [double]"99.51515143735945"
$ws.Range("B67").Value = [double]"1492.664422274971"
$ws.Range("C67").Value = [double]"746.3322131690232"
$ws.Range("D67").Value = [double]"-30.00020964141612"
$ws.Range("E67").Value = [double]"-30.00009391353698"
$ws.Range("F67").Value = [double]"-0.001134635891084315"
$ws.Range("G67").Value = [double]"1.726125258277022E-05"
$ws.Range("H67").Value = [double]"20.0981825890309"
$ws.Range("I67").Value = [double]"40.19670076572509"
$ws.Range("B68").Value = [double]"586.4879829573833"
$ws.Range("C68").Value = [double]"390.9920572034896"
$ws.Range("D68").Value = [double]"49.99969066422425"
$ws.Range("E68").Value = [double]"50.0022101445969"
$ws.Range("F68").Value = [double]"-0.006257545078988984"
$ws.Range("G68").Value = [double]"0.000475867418026269"
$ws.Range("H68").Value = [double]"85.25248643040224"
$ws.Range("B72").Value = [double]"19997.75723962632"
$ws.Range("C72").Value = [double]"19981.41654245573"
$ws.Range("D72").Value = [double]"582.6315257369653"
$ws.Range("E72").Value = [double]"582.1956903674263"
$ws.Range("F72").Value = [double]"66.998623982542"
$ws.Range("G72").Value = [double]"175.9199282817469"
$ws.Range("H72").Value = [double]"16.93186000817659"
$ws.Range("I72").Value = [double]"17.57337002406266"
$ws.Range("B73").Value = [double]"19997.75723962632"
$ws.Range("C73").Value = [double]"19991.44039087657"
$ws.Range("D73").Value = [double]"216.8396807166607"
$ws.Range("E73").Value = [double]"216.7774766492678"
$ws.Range("F73").Value = [double]"4.860091400329917"
$ws.Range("G73").Value = [double]"114.0867592606651"
$ws.Range("H73").Value = [double]"6.261888108244006"
$ws.Range("I73").Value = [double]"7.074625567814897"
$ws.Range("B74").Value = [double]"390.9920215025131"
$ws.Range("C74").Value = [double]"390.3614997407959"
$ws.Range("D74").Value = [double]"40.06457288426982"
$ws.Range("E74").Value = [double]"39.99999844378982"
$ws.Range("F74").Value = [double]"0.05249948831484329"
$ws.Range("G74").Value = [double]"-5.679085609244794E-08"
$ws.Range("H74").Value = [double]"59.16057266012713"
$ws.Range("I74").Value = [double]"59.16057266012713"
$ws.Range("B75").Value = [double]"19997.75723962632"
$ws.Range("C75").Value = [double]"19996.52758431912"
$ws.Range("D75").Value = [double]"93.57218479693776"
$ws.Range("E75").Value = [double]"93.56437262763318"
$ws.Range("F75").Value = [double]"-64.1677747988345"
$ws.Range("G75").Value = [double]"8.691112152841136"
$ws.Range("H75").Value = [double]"3.275636477606379"
$ws.Range("I75").Value = [double]"2.713100905171892"
$ws.Range("B76").Value = [double]"392.630725570652"
$ws.Range("C76").Value = [double]"383.6890695056338"
$ws.Range("D76").Value = [double]"204.1258983340154"
$ws.Range("E76").Value = [double]"199.9999770326304"
$ws.Range("F76").Value = [double]"100.2188166816432"
$ws.Range("G76").Value = [double]"96.86440911954156"
$ws.Range("H76").Value = [double]"334.3854571290636"
$ws.Range("I76").Value = [double]"334.3854571290636"
$ws.Range("B77").Value = [double]"392.630725570652"
$ws.Range("C77").Value = [double]"390.6949334566526"
$ws.Range("D77").Value = [double]"105.3023662202068"
$ws.Range("E77").Value = [double]"104.8396612734503"
$ws.Range("F77").Value = [double]"22.57881868669629"
$ws.Range("G77").Value = [double]"22.20263580315441"
$ws.Range("H77").Value = [double]"158.3698441106598"
$ws.Range("I77").Value = [double]"158.3698441106598"
$ws.Range("B78").Value = [double]"2009.81528168955"
$ws.Range("C78").Value = [double]"2029.411106824899"
$ws.Range("D78").Value = [double]"-278.0358415195804"
$ws.Range("E78").Value = [double]"-279.9999994790002"
$ws.Range("F78").Value = [double]"-167.278856110206"
$ws.Range("G78").Value = [double]"-169.7727224533312"
$ws.Range("H78").Value = [double]"93.21135910512457"
$ws.Range("I78").Value = [double]"93.15643607982241"
$ws.Range("B82").Value = [double]"586.0423178039709"
$ws.Range("C82").Value = [double]"586.4879829573833"
$ws.Range("D82").Value = [double]"-65.16277787202941"
$ws.Range("E82").Value = [double]"-65.21141776491942"
$ws.Range("F82").Value = [double]"-50.44338323266573"
$ws.Range("G82").Value = [double]"-50.4829205184717"
$ws.Range("H82").Value = [double]"140.6096317784082"
$ws.Range("I82").Value = [double]"140.609624112097"
